# Weekly update: insert a new price-report row (new row 39) that carries
# the previous week's "Haba" record forward, and refresh row 38 with the
# newest observation (new date + volume).
#
# Net effect vs. the original sheet:
#   - Row 38 keeps its identity but D38 (Fecha) becomes 44879 and
#     J38 (Volumen) becomes 120.
#   - A brand-new row is inserted right after it (becoming row 39) that
#     holds exactly what row 38 used to contain (D=44873, J=80, ...).
#   - Every row that used to be 38..63 is therefore shifted down to 39..64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 39 - this shifts old rows 39..63 down to
# 40..64, but leaves row 38 (and everything above it) untouched.
$ws.Rows.Item(39).Insert()

# Populate the newly-inserted row 39 with row 38's original data (a
# duplicate of the still-unmodified row 38, before we touch it below).
$ws.Range("A39").Value = 7
$ws.Range("B39").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C39").Value = "Ñuble"
$ws.Range("D39").Value = 44873
$ws.Range("D39").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("E39").Value = 16
$ws.Range("F39").Value = 100112026
$ws.Range("G39").Value = "Haba"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 80
$ws.Range("K39").Value = 6500
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = 6750
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Provincia de Diguillín"
$ws.Range("P39").Value = 270
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"

# Now refresh row 38 in place with this week's new observation: later
# date, bigger reported volume (prices/quality/origin unchanged).
$ws.Range("D38").Value = 44879
$ws.Range("J38").Value = 120
